$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.317.92'
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").Value = '1.860.45'
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7044'
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.19'
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07876'
$ws.Range("E8").Value = '  +1.78%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3048'
$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("E10").Value = '  +6.97%  '

$ws.Range("D11").Value = '2.022.69'
$ws.Range("E11").Value = '  +9.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08179'
$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.225'
$ws.Range("E13").Value = '  +1.16%  '

$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.64'
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("D16").Value = '29.346.43'
$ws.Range("E16").Value = '  +0.43%  '

$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '2.328.21'
$ws.Range("E17").Value = '  +10.45%  '

$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.841'
$ws.Range("E18").Value = '  +1.24%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007805'
$ws.Range("E19").Value = '  +0.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.27'
$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '238.42'
$ws.Range("E21").Value = '  +0.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.653'
$ws.Range("E24").Value = '  +2.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.89'
$ws.Range("E25").Value = '  +0.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.917'
$ws.Range("E26").Value = '  -1.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1431'
$ws.Range("E27").Value = '  -3.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.11'
$ws.Range("E28").Value = '  +0.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.906'
$ws.Range("E29").Value = '  -7.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.380'
$ws.Range("E30").Value = '  -3.49%  '

$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.328'
$ws.Range("E32").Value = '  -2.48%  '

$ws.Range("E33").Value = '  +0.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05194'
$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.179'
$ws.Range("E35").Value = '  +0.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7138'
$ws.Range("E36").Value = '  +0.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.010'
$ws.Range("E37").Value = '  +1.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.676'
$ws.Range("E38").Value = '  +0.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01854'
$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.689'
$ws.Range("E40").Value = '  -1.32%  '

$ws.Range("D41").Value = '1.176.47'
$ws.Range("E41").Value = '  +3.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9224'
$ws.Range("E42").Value = '  -1.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.021'
$ws.Range("E43").Value = '  +2.18%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.77'
$ws.Range("E44").Value = '  +1.50%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4279'
$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.199.97'
$ws.Range("E46").Value = '  +9.96%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9998'
$ws.Range("E47").Value = '  +0.03%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.70'
$ws.Range("E48").Value = '  -1.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5327'
$ws.Range("E49").Value = '  -2.23%  '

$ws.Range("E50").Value = '  -1.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.186'
$ws.Range("E51").Value = '  +0.11%  '
